$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF) - copy style from the existing
# header cell H1 (bold font, border, centered) before setting the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new I and J columns, keyed by row number.
$data = @{
    2  = @(6, 6)
    3  = @(6, 6)
    4  = @(5, 6)
    5  = @(7, 7)
    6  = @(7, 7)
    7  = @(9, 9)
    8  = @(8, 9)
    9  = @(6, 8)
    10 = @(6, 6)
    11 = @(8, 8)
    12 = @(4, 4)
    13 = @(7, 8)
    14 = @(3, 4)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(8, 8)
    18 = @(5, 5)
    19 = @(5, 6)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
